$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Statistics")
$ws.Range("D2").Value = 23.068616000000006
$ws.Range("E2").Value = 0.17849266249999118
$ws.Range("M2").Value = 249.49278027131
$ws.Range("N2").Value = -0.35826413287779246
$ws.Range("D3").Value = 19.262699750000003
$ws.Range("E3").Value = 7.305765810000011
$ws.Range("M3").Value = 82.86328482501447
$ws.Range("N3").Value = -0.6230701455942551
$ws.Range("D4").Value = 9.9607555
$ws.Range("E4").Value = 0.21198040000000695
$ws.Range("M4").Value = 111.49631916740344
$ws.Range("N4").Value = 1.00400953731679
$ws.Range("D5").Value = 19.964299749999977
$ws.Range("E5").Value = 7.4397893249999765
$ws.Range("M5").Value = 72.429970090252
$ws.Range("N5").Value = -1.0425820697849986
$ws.Range("D6").Value = 36.52823000000003
$ws.Range("E6").Value = -8.929370644999956
$ws.Range("M6").Value = 31.213017707021226
$ws.Range("N6").Value = 1.623694669381024
$ws.Range("D7").Value = 5.190070000000002
$ws.Range("E7").Value = 0.24951025500000146
$ws.Range("M7").Value = 41.726392798285815
$ws.Range("N7").Value = -1.899469376836315
$ws.Range("D8").Value = 13.219860749999997
$ws.Range("E8").Value = -0.41266167750000093
$ws.Range("M8").Value = 160.08442072665272
$ws.Range("N8").Value = 0.6927758803809922
$ws.Range("D9").Value = 5.878192749999999
$ws.Range("E9").Value = 0.00020594250000005587
$ws.Range("M9").Value = 87.0194978827494
$ws.Range("N9").Value = -1.1725074618526747
$ws.Range("D10").Value = 48.48797575
$ws.Range("E10").Value = 6.714489130000011
$ws.Range("M10").Value = 75.69028234346479
$ws.Range("N10").Value = 0.9670376983643223
$ws.Range("D11").Value = 19.225666000000004
$ws.Range("E11").Value = 2.2582373150000024
$ws.Range("M11").Value = 163.40211694872414
$ws.Range("N11").Value = -0.3699564803385442
$ws.Range("D12").Value = 4.1929412500000005
$ws.Range("E12").Value = -0.3727140924999963
$ws.Range("M12").Value = 26.03843802729898
$ws.Range("N12").Value = -0.1928227259734605
$ws.Range("D13").Value = 12.41103225
$ws.Range("E13").Value = 0.2272519699999993
$ws.Range("M13").Value = 153.6460497630287
$ws.Range("N13").Value = -0.5922196848961221
$ws.Range("D14").Value = 26.602334749999997
$ws.Range("E14").Value = 6.737270917500005
$ws.Range("M14").Value = 161.95816477345215
$ws.Range("N14").Value = 0.4893941677057114
$ws.Range("D15").Value = 243.99267450000002
$ws.Range("E15").Value = 21.60824731249997
$ws.Range("M15").Value = 1417.060735324658
$ws.Range("N15").Value = -1.4739801250050277

$ws = $wb.Worksheets.Item("Speeds")
$ws.Range("B2").Value = 10.815247012274597
$ws.Range("C2").Value = 0.08438564180090444
$ws.Range("H2").Value = 38.93488924418855
$ws.Range("I2").Value = 0.303788310483256
$ws.Range("B3").Value = 4.301748244039077
$ws.Range("C3").Value = 0.11314175948563998
$ws.Range("H3").Value = 15.486293678540678
$ws.Range("I3").Value = 0.4073103341483039
$ws.Range("B4").Value = 11.193560485186435
$ws.Range("C4").Value = 0.1373764100178965
$ws.Range("H4").Value = 40.29681774667117
$ws.Range("I4").Value = 0.4945550760644274
$ws.Range("B5").Value = 3.6279744843167907
$ws.Range("C5").Value = 0.16543973558617975
$ws.Range("H5").Value = 13.060708143540447
$ws.Range("I5").Value = 0.5955830481102471
$ws.Range("B6").Value = 0.8544902861984061
$ws.Range("C6").Value = 0.03585048006609994
$ws.Range("H6").Value = 3.076165030314262
$ws.Range("I6").Value = 0.1290617282379598
$ws.Range("B7").Value = 8.039658963807
$ws.Range("C7").Value = 0.05801435753794925
$ws.Range("H7").Value = 28.9427722697052
$ws.Range("I7").Value = 0.2088516871366173
$ws.Range("B8").Value = 12.109387818374165
$ws.Range("C8").Value = 0.1303152890754782
$ws.Range("H8").Value = 43.59379614614699
$ws.Range("I8").Value = 0.46913504067172157
$ws.Range("B9").Value = 14.803784357488007
$ws.Range("C9").Value = 0.2072659958984005
$ws.Range("H9").Value = 53.29362368695683
$ws.Range("I9").Value = 0.7461575852342418
$ws.Range("B10").Value = 1.5610113883433212
$ws.Range("C10").Value = 0.014437487048745134
$ws.Range("H10").Value = 5.6196409980359565
$ws.Range("I10").Value = 0.05197495337548248
$ws.Range("B11").Value = 8.499165487880841
$ws.Range("C11").Value = 0.1139337876947622
$ws.Range("H11").Value = 30.59699575637103
$ws.Range("I11").Value = 0.41016163570114395
$ws.Range("B12").Value = 6.2100650771815555
$ws.Range("C12").Value = 0.31579383682221224
$ws.Range("H12").Value = 22.356234277853602
$ws.Range("I12").Value = 1.136857812559964
$ws.Range("B13").Value = 12.37979619004122
$ws.Range("C13").Value = 0.20662683177290783
$ws.Range("H13").Value = 44.567266284148396
$ws.Range("I13").Value = 0.7438565943824682
$ws.Range("B14").Value = 6.088118441312831
$ws.Range("C14").Value = 0.05390599909432408
$ws.Range("H14").Value = 21.917226388726192
$ws.Range("I14").Value = 0.1940615967395667
$ws.Range("B15").Value = 5.807800329368733
$ws.Range("C15").Value = 0.041765108948216785
$ws.Range("H15").Value = 20.908081185727436
$ws.Range("I15").Value = 0.15035439221358043
